$d = $word.ActiveDocument
$find = $d.Content.Find
$find.Execute("Familiarity with HTML, CSS,", $true, $false, $false, $false, $false, $true, 1, $false, "Familiarity with Git, HTML, CSS,", 2)
